$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

$ws.Cells.Item($row, 1).Value = "'2025-09-22"
$ws.Cells.Item($row, 2).Value = 59.93000030517578
$ws.Cells.Item($row, 3).Value = 696.25
$ws.Cells.Item($row, 4).Value = 341.8500061035156

$ws.Cells.Item($row, 1).Style = "Normal"
